# Add a new slide (slide 9) that is a duplicate of slide 6 ("DiSCoVER: top
# drugs (cerebellar stem cell control)") and move it to the end of the deck,
# matching the new p:sldId id="264" r:id="rId15" entry appended to
# p:sldIdLst and the new ppt/slides/slide9.xml part added by the commit.

$p = $ppt.ActivePresentation

$source = $p.Slides.Item(6)
$newSlide = $source.Duplicate()

$newSlide.MoveTo($p.Slides.Count)
